$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: run a Find/Replace and fail loudly if the target text is missing,
# so problems surface immediately instead of silently producing a wrong doc.
# ---------------------------------------------------------------------------
function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute(
        $find, $true, $false, $false, $false, $false, $true, 1, $false,
        $replace, 2)
    if (-not $ok) {
        throw "Replace-Text failed to find: $find"
    }
}

# ---------------------------------------------------------------------------
# 1) Simple text replacements (do not change paragraph counts)
# ---------------------------------------------------------------------------

Replace-Text `
    "LLM inference efficiency and multi-agent systems." `
    "LLM inference efficiency, multi-agent systems, and AI safety research (activation probing, sandbagging detection)."

Replace-Text `
    "across 6 countries" `
    "across 7 countries"

Replace-Text `
    "serving 11 markets and 1200+ global users" `
    "serving 11 markets"

Replace-Text `
    "reduced model development time from 6 months to 1 week" `
    "reduced model development time from months to weeks"

Replace-Text `
    "Modernized MarTech infrastructure, driving 30% increase in customer acquisition" `
    "MarTech modernization - +30% customer acquisition"

Replace-Text `
    "Architected enterprise-scale data solutions for Fortune 500 clients across APAC, designing scalable platforms with measurable business impact." `
    "Architected enterprise-scale data solutions for Fortune 500 clients across APAC."

Replace-Text `
    "Engineered 5 high-performance data lakes processing 1.2 PB/hour, achieving 20% optimization" `
    "Data lakes processing 1.2 PB/hour for Fortune 500 clients across APAC"

Replace-Text `
    "Built real-time fraud detection systems, reducing false positives by 60% and saving `$XM annually" `
    "Real-time fraud detection systems - 60% reduction in false positives"

Replace-Text `
    "Various Companies" `
    "Microsoft, Truckaurbus (Founder), UTU"

Replace-Text `
    "Software Engineering, Architecture and Technical Consulting Roles" `
    "Software Engineering & Technical Leadership"

Replace-Text `
    "Progressively advanced through roles in software development, systems integration, and technical consulting within financial services and algorithmic trading domains." `
    "Progressive advancement through software engineering, entrepreneurship, and technical leadership across systems development, marketplace platforms, and payments infrastructure."

# ---------------------------------------------------------------------------
# 2) Insert three new bullet paragraphs after the "Progressive advancement..."
#    paragraph, describing Microsoft / Truckaurbus / UTU roles. We duplicate
#    an existing "ListParagraph"-styled bullet (paragraph 12) three times via
#    copy/paste so the new paragraphs inherit the correct style / numbering /
#    run formatting, then overwrite their text.
# ---------------------------------------------------------------------------

$bulletTemplateRange = $d.Paragraphs.Item(12).Range
$bulletTemplateRange.Copy() | Out-Null

# Locate the "Progressive advancement..." paragraph by scanning for its text.
# Note: Paragraph.Range.Text includes a trailing paragraph-mark character
# (chr 13), so it must be trimmed before comparing against plain text.
$progParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptxt = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($ptxt -eq "Progressive advancement through software engineering, entrepreneurship, and technical leadership across systems development, marketplace platforms, and payments infrastructure.") {
        $progParaIndex = $i
        break
    }
}
if ($progParaIndex -eq -1) {
    throw "Could not locate the 'Progressive advancement...' paragraph"
}

$insertPoint = $d.Paragraphs.Item($progParaIndex).Range.End
$r = $d.Range($insertPoint, $insertPoint)
$r.Paste() | Out-Null
$r = $d.Range($insertPoint, $insertPoint)
$r.Paste() | Out-Null
$r = $d.Range($insertPoint, $insertPoint)
$r.Paste() | Out-Null

$d.Paragraphs.Item($progParaIndex + 1).Range.Text = "Microsoft (2010-2014): Windows Kernel development (Windows 7/8, Server 2012 R2), Azure ML implementations, CDN architecture optimization"
$d.Paragraphs.Item($progParaIndex + 2).Range.Text = "Truckaurbus (2014-2016): Founded B2B commercial vehicle marketplace - 15 cities, 25+ OEM/bank partnerships"
$d.Paragraphs.Item($progParaIndex + 3).Range.Text = "UTU Singapore (2016-2017): Led maiden Thailand technical development; bank integration; payment/rewards systems for merchants"

# ---------------------------------------------------------------------------
# 3) Delete paragraphs that were removed outright. Deletions are performed by
#    locating each paragraph by its unique text and removing its whole range
#    (including the paragraph mark), from bottom to top so earlier indices
#    remain valid.
# ---------------------------------------------------------------------------

function Remove-ParagraphByText($text) {
    for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
        $ptxt = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
        if ($ptxt -eq $text) {
            $d.Paragraphs.Item($i).Range.Delete() | Out-Null
            return
        }
    }
    throw "Could not locate paragraph to remove: $text"
}

# Standard Chartered Bank bullet removed
Remove-ParagraphByText "Designed credit risk AI models integrating alternative data sources, improving accuracy by 15%"

# Think Big Analytics bullet removed
Remove-ParagraphByText "Designed enterprise architectures supporting global Fortune 500 clients across APAC"

# Entire "CatchMe - Intelligent Trust Engine" project block removed (3 paragraphs)
Remove-ParagraphByText "First of a kind, industry agnostic hybrid agentic AI decisioning system across Finance, Healthcare, Insurance, Cybersecurity, and Supply Chain. Uses adversarial debate protocols (prosecutor/defense/judge) to qualify events/anomalies and build audit trails for regulated environments."
Remove-ParagraphByText "Google Technical Disclosures - Pending (APLS & Cascade Routing)"
Remove-ParagraphByText "CatchMe - Intelligent Trust Engine (2025)"

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
